$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.400.84"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "2.192.48"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.70%  "

$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.67%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.05%  "

$ws.Range("E14").Value = "  -0.93%  "

$ws.Range("D15").Value = "2.516.06"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.872"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.15%  "

$ws.Range("E17").Value = "  -2.08%  "

$ws.Range("D18").Value = "2.184.01"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("D19").Value = "41.246.23"
$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("E21").Value = "  +2.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +20.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.37%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("E33").Value = "  -2.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.42%  "

$ws.Range("E35").Value = "  +1.62%  "

$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("E41").Value = "  -1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.204"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.98%  "

$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("E49").Value = "  +4.28%  "

$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
